$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 "I0" and J1 "IF" - copy formatting (style) from the
# existing H1 header cell so they share the same bold/bordered/centered style.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Fill data rows 2-21: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
